# Sky Swing - "Asset List - Code.docx" edits
#
# 1. Merge a couple of paragraphs whose sentence was split across several
#    <w:r> runs back into a single run (no visible/formatting change).
# 2. Cross out (strike-through) the Game-Manager background-tile spawning
#    bullets (now handled by a trigger) and the old mouse-aim/hook bullets
#    (superseded by the new trajectory+hook work), per the commit message.

$d = $word.ActiveDocument

function Set-ParagraphXml($Paragraph, $InnerXml) {
    # Replace the paragraph's content (everything except the trailing
    # paragraph mark) with a literal WordprocessingML fragment, so run
    # boundaries end up exactly as specified instead of however the
    # generic text-replace machinery decides to (re)split them.
    $r = $d.Range($Paragraph.Range.Start, $Paragraph.Range.End - 1)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

function Get-ParaText($p) {
    return $p.Range.Text.TrimEnd([char]13, [char]7)
}

# --- 1a. "Functions to " + "continue the level" + ", restart the level, or
#          return to main menu" -> one run (tab stays a separate <w:tab/>).
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "`tFunctions to continue the level, restart the level, or return to main menu") {
        Set-ParagraphXml $p '<w:r><w:tab/><w:t>Functions to continue the level, restart the level, or return to main menu</w:t></w:r>'
        break
    }
}

# --- 1b. "Boost forward in direction of movement with M2" + " if you have
#          boost left" -> one run.
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "Boost forward in direction of movement with M2 if you have boost left") {
        Set-ParagraphXml $p '<w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Boost forward in direction of movement with M2 if you have boost left</w:t></w:r>'
        break
    }
}

# --- 1c. "M1 attaches hook" + " onto building" -> one run (strike-through
#          applied afterwards along with the other crossed-out bullets).
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "M1 attaches hook onto building") {
        Set-ParagraphXml $p '<w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>M1 attaches hook onto building</w:t></w:r>'
        break
    }
}

# --- 2. Strike-through the bullets that are now crossed out.
$strikeTexts = @(
    "Manage the different background tiles, spawning them as the player gets past certain boundaries to make endless level",
    "Maybe Move the background with the player movement",
    "Mouse pointer to aim",
    "M1 attaches hook onto building",
    "Once hook is attached, M1 and hold shortens rope"
)

foreach ($p in $d.Paragraphs) {
    if ($strikeTexts -contains (Get-ParaText $p)) {
        $p.Range.Font.StrikeThrough = 1
    }
}
